$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update the tracked google ad-click URL ---
$ws.Range("B2").Value = "https://www.google.com/aclk?sa=l&ai=DChcSEwj8tZ3p7fyCAxWnZJEFHTGjCBUYABAAGgJscg&ase=2&gclid=EAIaIQobChMI_LWd6e38ggMVp2SRBR0xowgVEBAYASAAEgIJHfD_BwE&sig=AOD64_1fMg5w2jnuxVVBjyMpP9pU9MCXmA&adurl&ctype=99&nis=4&ved=2ahUKEwiL5I_p7fyCAxUUGwYAHftOBE4Q8PwKegQIABAo"

# --- Insert a new row 3 for "makeagency" (pushes old row3 "google" down to row4, etc.) ---
$ws.Rows.Item(3).EntireRow.Insert()
$ws.Range("A3").Value = "makeagency"
$ws.Range("B3").Value = "http://makeagency.co.uk/"
$ws.Range("C3").Value = "London"
$ws.Range("D3").Value = "UK"
$ws.Range("E3").Value = "hi@makeagency.co.uk"

# --- Row 5 (was "bartleboglehegarty", now shifted) becomes "londonmarketingcompany" ---
$ws.Range("A5").Value = "londonmarketingcompany"
$ws.Range("B5").Value = "http://www.londonmarketingcompany.co.uk/"
$ws.Range("C5").Value = "London"
$ws.Range("D5").Value = "UK"
$ws.Range("E5").Value = "info@bathmarketingcompany.uixweb.dev"

# --- Insert a new row 9 for "unrvld" (pushes old row9 "vccp" down to row10, etc.) ---
$ws.Rows.Item(9).EntireRow.Insert()
$ws.Range("A9").Value = "unrvld"
$ws.Range("B9").Value = "https://unrvld.com/"
$ws.Range("C9").Value = "London"
$ws.Range("D9").Value = "UK"
$ws.Range("E9").Value = "hello@unrvld.com, marketing@unrvld.com, vacancies@unrvld.com"
